$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue 'D2' '62.067.48'
Set-TextValue 'E2' '  -0.25%  '

Set-TextValue 'D3' '3.424.08'
Set-TextValue 'E3' '  -0.18%  '

Set-TextValue 'D4' '1.00'
Set-TextValue 'E4' '  +0.09%  '

Set-TextValue 'D5' '410.97'
Set-TextValue 'E5' '  +0.91%  '

Set-TextValue 'D6' '130.01'
Set-TextValue 'E6' '  -3.05%  '

Set-TextValue 'D7' '0.645'
Set-TextValue 'E7' '  +9.18%  '

Set-TextValue 'E8' '  -0.05%  '

Set-TextValue 'E9' '  +7.95%  '

Set-TextValue 'D10' '0.143'
Set-TextValue 'E10' '  +16.02%  '

Set-TextValue 'D11' '42.80'
Set-TextValue 'E11' '  +1.25%  '

Set-TextValue 'D12' '0.0000221'
Set-TextValue 'E12' '  +68.22%  '

Set-TextValue 'B13' 'Polkadot'
Set-TextValue 'C13' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D13' '9.17'
Set-TextValue 'E13' '  +8.60%  '

Set-TextValue 'B14' 'TRON'
Set-TextValue 'C14' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue 'D14' '0.141'
Set-TextValue 'E14' '  -0.36%  '

Set-TextValue 'D15' '3.957.54'
Set-TextValue 'E15' '  -0.39%  '

Set-TextValue 'D16' '20.98'
Set-TextValue 'E16' '  +5.11%  '

Set-TextValue 'D17' '3.429.55'
Set-TextValue 'E17' '  +0.37%  '

Set-TextValue 'D18' '12.12'
Set-TextValue 'E18' '  +7.21%  '

Set-TextValue 'D19' '1.08'
Set-TextValue 'E19' '  +5.86%  '

Set-TextValue 'D20' '62.045.17'
Set-TextValue 'E20' '  -0.26%  '

Set-TextValue 'D21' '441.30'
Set-TextValue 'E21' '  +40.45%  '

Set-TextValue 'D22' '91.06'
Set-TextValue 'E22' '  +7.85%  '

Set-TextValue 'D23' '3.18'

Set-TextValue 'D24' '13.12'
Set-TextValue 'E24' '  +1.44%  '

Set-TextValue 'D25' '3.26'
Set-TextValue 'E25' '  +3.43%  '

Set-TextValue 'D26' '33.81'
Set-TextValue 'E26' '  +13.75%  '

Set-TextValue 'D27' '8.86'
Set-TextValue 'E27' '  +8.38%  '

Set-TextValue 'E28' '  -0.11%  '

Set-TextValue 'D29' '7.69'
Set-TextValue 'E29' '  +1.63%  '

Set-TextValue 'D30' '2.70'
Set-TextValue 'E30' '  -3.44%  '

Set-TextValue 'D31' '12.06'
Set-TextValue 'E31' '  +5.87%  '

Set-TextValue 'D32' '0.116'
Set-TextValue 'E32' '  -0.23%  '

Set-TextValue 'E33' '  -2.64%  '

Set-TextValue 'D34' '42.87'
Set-TextValue 'E34' '  +0.98%  '

Set-TextValue 'E35' '  -0.07%  '

Set-TextValue 'D36' '0.0502'
Set-TextValue 'E36' '  +3.19%  '

Set-TextValue 'D37' '53.88'
Set-TextValue 'E37' '  +4.64%  '

Set-TextValue 'D38' '0.998'
Set-TextValue 'E38' '  +0.06%  '

Set-TextValue 'D39' '0.136'
Set-TextValue 'E39' '  +8.18%  '

Set-TextValue 'D40' '3.39'
Set-TextValue 'E40' '  -0.49%  '

Set-TextValue 'E41' '  -1.02%  '

Set-TextValue 'D42' '0.317'
Set-TextValue 'E42' '  +1.76%  '

Set-TextValue 'D43' '141.74'
Set-TextValue 'E43' '  +2.36%  '

Set-TextValue 'D44' '1.99'
Set-TextValue 'E44' '  +0.03%  '

Set-TextValue 'D45' '4.15'
Set-TextValue 'E45' '  +2.92%  '

Set-TextValue 'E46' '  +7.86%  '

Set-TextValue 'E47' '  -0.90%  '

Set-TextValue 'D48' '22.21'
Set-TextValue 'E48' '  +4.20%  '

Set-TextValue 'D49' '3.769.38'
Set-TextValue 'E49' '  +0.06%  '

Set-TextValue 'B50' 'Maker'
Set-TextValue 'C50' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D50' '2.116.86'
Set-TextValue 'E50' '  -0.21%  '

Set-TextValue 'B51' 'BitcoinSV'
Set-TextValue 'C51' 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue 'D51' '105.51'
Set-TextValue 'E51' '  +26.32%  '
